$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bValues = @(4.47, 4.46, 4.35, 4.28, 4.44, 4.24, 4.35, 4.3, 4.37, 4.67, 4.29, 4.39)
$cValues = @(5.12, 4.85, 4.91, 5.03, 5.23, 5, 4.95, 4.87, 4.94, 5.02, 4.97, 4.87)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}
